$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is updated for every
# data row (rows 2-257) from 45175 to 45177.
$ws.Range("C2:C257").Value = 45177
